# Autenticacion.xlsx data update
# - Replace obsolete test users/error messages with the new ones used by the
#   updated Android/iOS authentication steps.
# - "Datos" (sheet1) is the only sheet whose cell content actually changes;
#   "Prueba" and "Listas" keep the same values (only internal shared-string
#   bookkeeping shifts, which Excel/this engine handles on its own).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# Row 2 ("ID"=1): usuario autotest30 -> pruebauser01
$ws.Range("D2").Value = "pruebauser01"

# Row 6 ("ID"=5): usuario testing10 -> pruebauser01
$ws.Range("D6").Value = "pruebauser01"

# Row 7 ("ID"=6): usuario userrobot9 -> bloqueotarjetas01,
# resultadoEsperado "Clave bloqueada" -> "La clave que usas en el cajero está bloqueada"
$ws.Range("D7").Value = "bloqueotarjetas01"
$ws.Range("J7").Value = "La clave que usas en el cajero está bloqueada"

# Row 6 resultadoEsperado "CLAVE NO VALIDA" -> "clave inválida"
$ws.Range("J6").Value = "clave inválida"

# Restore the active selection to the cell the author left selected (E6)
$ws.Activate()
$ws.Range("E6").Select() | Out-Null
